$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (appears on Overview!E2:F3 and on the zh-cn / de-de Status column C2:C3)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Handback report details for the "zh-cn" sheet (rows 2 & 3):
#      I = Latest Target File    -> "a.md" (hyperlinked, like column A)
#      J = Latest Handback File  -> generated xlf file name
#      K = Latest Handback DateTime
# ---------------------------------------------------------------------------
$zhHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhHandbackDate = "2016-09-04 02:40:54"
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/816db466cb003b699058d14ac759f5769d88a66e/e2e/a.md"

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsZh.Range("J2").Value = $zhHandbackFile
$wsZh.Range("K2").Value = $zhHandbackDate

$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsZh.Range("J3").Value = $zhHandbackFile
$wsZh.Range("K3").Value = $zhHandbackDate

# ---------------------------------------------------------------------------
# 3. Handback report details for the "de-de" sheet (rows 2 & 3):
# ---------------------------------------------------------------------------
$deHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deHandbackDate = "2016-09-04 02:41:04"

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsDe.Range("J2").Value = $deHandbackFile
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aMdUrl, [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsDe.Range("J3").Value = $deHandbackFile
$wsDe.Range("K3").Value = $deHandbackDate

# ---------------------------------------------------------------------------
# 4. Widen columns that now hold the longer status text / file names.
#    (ColumnWidth is quantised by the host; these values land on the exact
#    pixel-rounded widths used elsewhere in this workbook.)
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.15   # E: zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = 29.15   # F: de-de status column

$wsZh.Columns.Item(3).ColumnWidth = 29.15     # C: Status
$wsZh.Columns.Item(10).ColumnWidth = 39.15    # J: Latest Handback File

$wsDe.Columns.Item(3).ColumnWidth = 29.15     # C: Status
$wsDe.Columns.Item(10).ColumnWidth = 39.15    # J: Latest Handback File
